$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRYNF")

# Row 17 - Total Operating Expenses
$ws.Range("J17").Value = 2800

# Row 18 - Operating Income or Loss
$ws.Range("J18").Value = -900

# Row 23 - Income Before Tax
$ws.Range("J23").Value = -900

# Row 60 - Total Current Liabilities
$ws.Range("J60").Value = 800

# Row 66 - Total Liabilities
$ws.Range("J66").Value = 800

# Row 72 - Retained Earnings
$ws.Range("D72").Value = -5000
$ws.Range("E72").Value = -4800
$ws.Range("F72").Value = -4600
$ws.Range("G72").Value = -4600
$ws.Range("H72").Value = -4500
$ws.Range("I72").Value = -4500
$ws.Range("J72").Value = -4200

# Row 76 - Total Stockholder Equity
$ws.Range("D76").Value = 500

# Row 91 - Capital Expenditures
$ws.Range("D91").Value = "NA"
$ws.Range("E91").Value = "NA"

# Row 100 - Other Cash Flows from Financing Activities
$ws.Range("E100").Value = 100
$ws.Range("J100").Value = 200

# Row 102 - Change In Cash and Cash Equivalents
$ws.Range("D102").Value = 400
